$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 6 with the latest quotations (2025-09-10)
$ws.Range("A6").Value = 45910
$ws.Range("A6").Style = $ws.Range("A2").Style
$ws.Range("A6").NumberFormat = $ws.Range("A2").NumberFormat

$ws.Range("B6").Value = "20,9112"
$ws.Range("C6").Value = "14,7555"
$ws.Range("D6").Value = "14,8193"
$ws.Range("E6").Value = "14,8193"
